$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.979.02"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.362.65"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'310.25"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'107.51"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'40.55"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'8.39"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'0.968"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "2.723.78"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'15.17"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "2.367.00"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "44.959.97"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "'14.25"
$ws.Range("E19").Value = "  +9.95%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("D21").Value = "'0.0000105"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'73.05"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'3.46"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'258.71"
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'11.06"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "'7.15"
$ws.Range("E28").Value = "  -5.20%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'22.28"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "'0.0959"
$ws.Range("E31").Value = "  +7.68%  "
$ws.Range("D32").Value = "'37.03"
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "'168.20"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  +5.04%  "
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").Value = "'0.115"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "'4.70"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'3.89"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").Value = "'2.90"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'0.0351"
$ws.Range("E40").Value = "  -3.62%  "
$ws.Range("D41").Value = "'1.73"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "'99.95"
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'69.34"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.865.30"
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("D45").Value = "'0.227"
$ws.Range("E45").Value = "  -5.35%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'12.69"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'5.62"
$ws.Range("E48").Value = "  +8.76%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'80.35"
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("D50").Value = "'111.19"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").Value = "'9.13"
$ws.Range("E51").Value = "  +1.66%  "

# Reset style on cells where we used a leading apostrophe to force text,
# so the quotePrefix flag introduced by the apostrophe does not linger.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
